$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: insert a new row at 15, shifting the existing row 15
#            ("Segunda", Castle Brite) down to row 16 ---
$ws.Rows("15").Insert()

# --- Step 2: fill the freshly inserted row 15 with the data that used to
#            live in row 14 (a duplicate "Primera" / Castle Brite record,
#            unchanged from before the edit) ---
$ws.Range("A15").Value = 11
$ws.Range("B15").Value = "Vega Monumental Concepción"
$ws.Range("C15").Value = "Bíobío"
$ws.Range("D15").Value = 44545
$ws.Range("E15").Value = 8
$ws.Range("F15").Value = "Fruta"
$ws.Range("G15").Value = 100103
$ws.Range("H15").Value = "Frutos de hueso (carozo)"
$ws.Range("I15").Value = 100103003
$ws.Range("J15").Value = "Damasco"
$ws.Range("K15").Value = "Castle Brite"
$ws.Range("L15").Value = "Primera"
$ws.Range("M15").Value = 100
$ws.Range("N15").Value = 18000
$ws.Range("O15").Value = 19000
$ws.Range("P15").Value = 18500
$ws.Range("Q15").Value = "$/caja 15 kilos"
$ws.Range("R15").Value = "Región de O'Higgins"
$ws.Range("S15").Value = 1233
$ws.Range("T15").Value = 15

# --- Step 3: update row 14 in place with the new week's record ---
$ws.Range("D14").Value = 44918
$ws.Range("K14").Value = "Dina"
$ws.Range("N14").Value = 17000
$ws.Range("O14").Value = 18000
$ws.Range("P14").Value = 17500
$ws.Range("Q14").Value = "$/caja 18 kilos"
$ws.Range("S14").Value = 972
$ws.Range("T14").Value = 18
